# Fix "Recorded By" (column G) entries so that the "System" / "system"
# token no longer sorts first in the comma-separated list of recorders.
# Only the two exact source strings that appear in the workbook are
# touched, matching the upstream diff precisely:
#   "System, dnasr281@gmail.com"           -> "dnasr281@gmail.com, System"
#   "system, System, backup@backdoor.com"  -> "System, system, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colG = 7

$replacements = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "system, System, backup@backdoor.com" = "System, system, backup@backdoor.com";
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $rawValue = $cell.Value2

    if ($rawValue -eq $null) { continue }

    $value = [string]$rawValue
    if ([string]::IsNullOrEmpty($value)) { continue }

    if ($replacements.ContainsKey($value)) {
        $cell.Value = $replacements[$value]
    }
}
